$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.339.73"
$ws.Range("E2").Value = "  +3.34%  "

# Row 3
$ws.Range("D3").Value = "2.250.83"
$ws.Range("E3").Value = "  +1.88%  "

# Row 4
$ws.Range("E4").Value = "  -0.19%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.67"
$ws.Range("E5").Value = "  +2.89%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "91.72"
$ws.Range("E6").Value = "  +4.46%  "

# Row 7
$ws.Range("E7").Value = "  +2.12%  "

# Row 8
$ws.Range("E8").Value = "  -0.12%  "

# Row 9
$ws.Range("E9").Value = "  +2.25%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "53.75"
$ws.Range("E10").Value = "  +7.75%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "32.10"
$ws.Range("E11").Value = "  +6.84%  "

# Row 12
$ws.Range("E12").Value = "  +2.16%  "

# Row 13
$ws.Range("E13").Value = "  +2.96%  "

# Row 14
$ws.Range("E14").Value = "  +2.20%  "

# Row 15
$ws.Range("D15").Value = "2.596.78"
$ws.Range("E15").Value = "  +1.69%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.14"
$ws.Range("E16").Value = "  +2.98%  "

# Row 17
$ws.Range("D17").Value = "2.251.49"
$ws.Range("E17").Value = "  +3.50%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.749"
$ws.Range("E18").Value = "  +3.31%  "

# Row 19
$ws.Range("D19").Value = "41.216.53"
$ws.Range("E19").Value = "  +3.18%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.19"
$ws.Range("E20").Value = "  +8.36%  "

# Row 21
$ws.Range("D21").Value = "0.0₃0903"
$ws.Range("E21").Value = "  +2.30%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.87"
$ws.Range("E22").Value = "  +1.81%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "66.72"
$ws.Range("E23").Value = "  +2.33%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "240.55"
$ws.Range("E24").Value = "  +1.36%  "

# Row 25
$ws.Range("E25").Value = "  +3.95%  "

# Row 26
$ws.Range("E26").Value = "  -0.15%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.87"
$ws.Range("E27").Value = "  +3.21%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "23.77"
$ws.Range("E28").Value = "  +5.85%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.18"
$ws.Range("E29").Value = "  +0.98%  "

# Row 30
$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.63"
$ws.Range("E30").Value = "  +5.29%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "158.55"
$ws.Range("E31").Value = "  +1.35%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "33.56"
$ws.Range("E32").Value = "  +7.25%  "

# Row 33
$ws.Range("E33").Value = "  -0.10%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.19"
$ws.Range("E34").Value = "  +6.41%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0735"
$ws.Range("E35").Value = "  +3.76%  "

# Row 36
$ws.Range("E36").Value = "  +7.51%  "

# Row 37
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.36"
$ws.Range("E37").Value = "  +1.16%  "

# Row 38
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "16.70"
$ws.Range("E38").Value = "  +8.69%  "

# Row 39
$ws.Range("E39").Value = "  +2.97%  "

# Row 40
$ws.Range("E40").Value = "  +5.77%  "

# Row 41
$ws.Range("E41").Value = "  +6.81%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.95"
$ws.Range("E42").Value = "  +5.88%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.27"
$ws.Range("E43").Value = "  +16.55%  "

# Row 44
$ws.Range("D44").Value = "2.068.11"
$ws.Range("E44").Value = "  -2.49%  "

# Row 45
$ws.Range("E45").Value = "  +3.22%  "

# Row 46
$ws.Range("E46").Value = "  +6.09%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.96"
$ws.Range("E47").Value = "  +11.69%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.05"
$ws.Range("E48").Value = "  -2.56%  "

# Row 49
$ws.Range("D49").Value = "2.468.30"
$ws.Range("E49").Value = "  +1.91%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.51"
$ws.Range("E50").Value = "  +1.74%  "

# Row 51
$ws.Range("E51").Value = "  +3.67%  "
